# Rewrites "game rules and output format.docx" from the Rock-Paper-Scissors
# rules text to the "silence / betrayal" rules text, per the target diff.
#
# Host quirks worked around here:
#  1. Setting `.Text` directly on the Range object returned by the
#     `Paragraph.Range` property *inserts* instead of *replacing* the
#     paragraph's text. Re-wrapping the same (Start, End) pair with
#     `$d.Range(start, end)` gives a Range whose `.Text` setter performs a
#     proper replace, so paragraph rewrites go through that instead.
#  2. For a paragraph that currently holds no real text (only the paragraph
#     mark, i.e. End - Start <= 1), wrapping the *full* (Start, End) span
#     and assigning `.Text` corrupts neighbouring text (it only partially
#     writes and the remainder leaks into the next paragraph). Using a
#     zero-width Range collapsed at Start instead makes `.Text =` behave
#     as a clean insert, which is what's needed for an originally-empty
#     paragraph anyway.
#
# Strategy: walk the document's paragraphs from LAST to FIRST so that the one
# paragraph insertion we need (a brand-new "2. silence and betrayal ..." rule)
# never shifts the positional index of a paragraph we still have to touch.

$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $s = $p.Range.Start
    $e = $p.Range.End
    if ($e - $s -le 1) {
        $r = $d.Range($s, $s)
    } else {
        $r = $d.Range($s, $e)
    }
    $r.Text = $newText
    return $r
}

# ---------------------------------------------------------------------
# Paragraph 9 (last rules-ish paragraph): the old "paper/stone/scissors"
# sentence becomes a blank paragraph.
# ---------------------------------------------------------------------
Set-ParaText 9 "" | Out-Null

# ---------------------------------------------------------------------
# Paragraph 8: "You can only output one word" -> the new "or" sentence.
# ---------------------------------------------------------------------
Set-ParaText 8 "‘silence’  or  ‘betrayal’   " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 7: "Output rules and formats" -> "You can only output one word "
# ---------------------------------------------------------------------
Set-ParaText 7 "You can only output one word " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 6: the blank paragraph -> "Output rules and formats "
# ---------------------------------------------------------------------
Set-ParaText 6 "Output rules and formats " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 5: "4.If both players choose the same move, the game is a draw."
#   -> a single space " "
# ---------------------------------------------------------------------
Set-ParaText 5 " " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 4: "3.Stone beats Scissors" -> "4.try to get as many points as you can  "
# ---------------------------------------------------------------------
Set-ParaText 4 "4.try to get as many points as you can  " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 3: "2.Paper beats Stone" -> "3.betrayal and betrayal, get 18 points  "
# ---------------------------------------------------------------------
Set-ParaText 3 "3.betrayal and betrayal, get 18 points  " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 2: "1.Scissors beats Paper" -> "1.silence and silence ,  get one point "
# then a brand-new paragraph is inserted right after it for rule 2.
# ---------------------------------------------------------------------
$r2 = Set-ParaText 2 "1.silence and silence ,  get one point "
$r2.InsertParagraphAfter()
Set-ParaText 3 "2.silence and betrayal,  get 20 points,betrayal get 0 pooint  " | Out-Null

# ---------------------------------------------------------------------
# Paragraph 1: add a trailing space.
# ---------------------------------------------------------------------
Set-ParaText 1 "The rules of the game are as follows: " | Out-Null
